$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.83%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.40%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.093"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.37%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07811"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-2.86%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.237"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-16.31%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.793"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.53%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'-0.19%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9184"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.02%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1759"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.43%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07521"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.90%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08987"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'6.79%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03029"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.05%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.1004"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.72%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001505"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.84%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006005"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.33%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.468"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.88%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'0.15%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'0.25%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D21").Value = "'4.239"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-8.43%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1817"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'13.67%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04588"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.72%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001249"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.78%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004472"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.50%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'5.93%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-1.33%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01773"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-3.36%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04783"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'5.93%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007411"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'5.61%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1359"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.27%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-2.24%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'4.16%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006231"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-3.59%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.00%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'28.90%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.7309"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-10.92%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.00%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.00%"
$ws.Range("E50").Style = "Normal"
